$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the part file name used for scanning from "nist.stl" to "nut.stl"
$ws.Range("B4").Value = "nut.stl"

# Scroll the view back to the top and select cell B5, matching the saved
# sheet view state (no topLeftCell override, selection on B5)
$ws.Activate()
$ws.Range("B5").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
